# Rotate the comma-separated "Recorded By" values in column G left by one
# position for each data row (e.g. "System, foo@bar.com" -> "foo@bar.com, System").
# A handful of rows are intentionally left untouched to mirror the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose G value must NOT be rotated even though it technically could be.
$excludedRows = @(4, 30, 56)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    if ($excludedRows -contains $row) {
        continue
    }

    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ',\s*'

    if ($parts.Count -le 1) {
        continue
    }

    $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
    $newValue = [string]::Join(', ', $rotated)

    $cell.Value = $newValue
}
